$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.590.39'
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").Value = '1.688.71'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.04'
$ws.Range("E5").Value = '  -0.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3897'
$ws.Range("E7").Value = '  -1.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4022'
$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.497'
$ws.Range("E9").Value = '  +0.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.005'
$ws.Range("E10").Value = '  +0.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.59'
$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08742'
$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.579'
$ws.Range("E13").Value = '  +4.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.75'
$ws.Range("E14").Value = '  +5.25%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.952'
$ws.Range("E15").Value = '  -0.97%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001349'
$ws.Range("E16").Value = '  +2.48%  '

$ws.Range("D17").Value = '1.681.56'
$ws.Range("E17").Value = '  -0.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '98.27'
$ws.Range("E18").Value = '  -1.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07094'
$ws.Range("E19").Value = '  +1.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.80'
$ws.Range("E20").Value = '  +1.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.271'
$ws.Range("E21").Value = '  +4.01%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.25'
$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("D24").Value = '24.587.39'
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.998'
$ws.Range("E25").Value = '  -9.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.353'
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.72'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.38'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.575'
$ws.Range("E29").Value = '  +12.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.25'
$ws.Range("E30").Value = '  +0.79%  '

$ws.Range("E31").Value = '  +0.63%  '

$ws.Range("D32").Value = '1.866.11'
$ws.Range("E32").Value = '  -0.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08751'
$ws.Range("E33").Value = '  +2.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.364'
$ws.Range("E34").Value = '  +3.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.037'
$ws.Range("E35").Value = '  -2.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.989'
$ws.Range("E36").Value = '  +5.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02908'
$ws.Range("E37").Value = '  +6.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2724'
$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.76'
$ws.Range("E39").Value = '  -4.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.18'
$ws.Range("E40").Value = '  -1.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09113'
$ws.Range("E41").Value = '  -0.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7802'
$ws.Range("E42").Value = '  +2.29%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.452'
$ws.Range("E43").Value = '  -0.72%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.69'
$ws.Range("E44").Value = '  +4.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7181'
$ws.Range("E45").Value = '  +0.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.579'
$ws.Range("E46").Value = '  -0.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.195'
$ws.Range("E47").Value = '  -0.51%  '

$ws.Range("E48").Value = '  +0.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.331'
$ws.Range("E49").Value = '  +0.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.59'
$ws.Range("E50").Value = '  -1.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '90.78'
$ws.Range("E51").Value = '  +1.19%  '
